$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Workbook window size (bookViews/workbookView)
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Width = 28800
$wb.Windows.Item(1).Height = 12330

# ---------------------------------------------------------------------------
# 2. Sheet zoom + selection
# ---------------------------------------------------------------------------
$null = $ws.Select()
$excel.ActiveWindow.Zoom = 130
$null = $ws.Range("K15").Select()

# ---------------------------------------------------------------------------
# 3. Remove the "192.1.30.1" row (old row 8) -- everything below it in the
#    C/D (router label) & E/F/G (ip / mask / soho) block shifts up one row.
#    Handled here with direct value writes (merged-cell aware) rather than
#    a native row/range Delete-shift, which this block's merges make
#    unreliable.
# ---------------------------------------------------------------------------

# C/D label column (merged pairs) - "4. szint router" moves row9 -> row8,
# "3. szint router" moves row12 -> row11
$ws.Range("C8").Value2 = $ws.Range("C9").Value2
$ws.Range("C9:D9").UnMerge()
$ws.Range("C9:D9").Clear()

$ws.Range("C2:D2").Copy()
$ws.Range("C11:D11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C11:D11").Merge()
$ws.Range("C11").Value2 = $ws.Range("C12").Value2
$ws.Range("C12:D12").UnMerge()
$ws.Range("C12:D12").Clear()

# E/F (+G) data column - shift every row 9..14 up into 8..13
$ws.Range("E8").Value2 = $ws.Range("E9").Value2
$ws.Range("F8").Value2 = $ws.Range("F9").Value2

$ws.Range("E9").Value2 = $ws.Range("E10").Value2
$ws.Range("F9").Value2 = $ws.Range("F10").Value2

$ws.Range("E10").Value2 = $ws.Range("E11").Value2
$ws.Range("F10").Value2 = $ws.Range("F11").Value2
$ws.Range("G10").Value2 = $ws.Range("G11").Value2
$ws.Range("G11").Value2 = ""

$ws.Range("E11").Value2 = $ws.Range("E12").Value2
$ws.Range("F11").Value2 = $ws.Range("F12").Value2

$ws.Range("E12").Value2 = $ws.Range("E13").Value2
$ws.Range("F12").Value2 = $ws.Range("F13").Value2

$ws.Range("E13").Value2 = $ws.Range("E14").Value2
$ws.Range("F13").Value2 = $ws.Range("F14").Value2
$ws.Range("G13").Value2 = $ws.Range("G14").Value2

$ws.Range("E14").Value2 = ""
$ws.Range("F14").Value2 = ""
$ws.Range("G14").Value2 = ""

# ---------------------------------------------------------------------------
# 4. Remove the "Vezetoseg_wifi" entry (old I12). The I column labels
#    (wifi / Wireless router0 / Wireless router1) stay put; only the
#    paired K/L subnet values shift down one row, and the last one
#    ("192.168.40.0") is dropped.
# ---------------------------------------------------------------------------
$ws.Range("I12").Value2 = ""

$ws.Range("K15").Value2 = $ws.Range("K14").Value2
$ws.Range("L15").Value2 = $ws.Range("L14").Value2

$ws.Range("K14").Value2 = $ws.Range("K13").Value2
$ws.Range("L14").Value2 = $ws.Range("L13").Value2

$ws.Range("K13").Value2 = $ws.Range("K12").Value2
$ws.Range("L13").Value2 = $ws.Range("L12").Value2

$ws.Range("K12").Value2 = ""
$ws.Range("L12").Value2 = ""

# ---------------------------------------------------------------------------
# 5. New blank styled/merged placeholder rows further down the sheet
#    (C19:D19, C20:D20, C21:D21, C24:D24)
# ---------------------------------------------------------------------------
$ws.Range("C2:D2").Copy()
foreach ($r in 19, 20, 21, 24) {
    $rng = $ws.Range("C" + $r + ":D" + $r)
    $rng.PasteSpecial(-4122)   # xlPasteFormats
    $rng.Merge()
}

$excel.CutCopyMode = $false
